$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New data row (ID=2, Creditor-Name=name2, Amount=10000, Date=1397:6:15).
# ID and Amount look like plain numbers, but the source data (like row 2)
# stores them as text, so force text storage for those two cells, then
# restore the default "Normal" style so the cell itself carries no
# explicit style index (matching the rest of the data rows).
$ws.Range("A3").NumberFormat = "@"
$ws.Range("A3").Value = "2"
$ws.Range("A3").Style = "Normal"

$ws.Range("B3").Value = "name2"

$ws.Range("C3").NumberFormat = "@"
$ws.Range("C3").Value = "10000"
$ws.Range("C3").Style = "Normal"

$ws.Range("D3").Value = "1397:6:15"
